$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns E-H are identical across rows 2-7
$ws.Range("E2:E7").Value = 3
$ws.Range("F2:F7").Value = 1
$ws.Range("G2:G7").Value = 0.8785453333333333
$ws.Range("H2:H7").Value = 2.635636

# Row 2
$ws.Range("M2").Value = 4.260872666666667
$ws.Range("N2").Value = 12.782618
$ws.Range("O2").Value = 0.09064705929364961
$ws.Range("P2").Value = 0.09064705929364959
$ws.Range("Q2").Value = 3.743369797227555
$ws.Range("R2").Value = 33.690328175048
$ws.Range("S2").Value = 0.09064705929364961
$ws.Range("T2").Value = 0.09064705929364959

# Row 3 (M3, N3 unchanged)
$ws.Range("O3").Value = 0.4366505728284585
$ws.Range("P3").Value = 0.4366505728284584
$ws.Range("Q3").Value = 18.031964621964
$ws.Range("R3").Value = 162.287681597676
$ws.Range("S3").Value = 0.4366505728284585
$ws.Range("T3").Value = 0.4366505728284584

# Row 4
$ws.Range("M4").Value = 8.931090666666666
$ws.Range("N4").Value = 26.793272
$ws.Range("O4").Value = 0.190002651698962
$ws.Range("P4").Value = 0.1900026516989619
$ws.Range("Q4").Value = 7.846368026776888
$ws.Range("R4").Value = 70.61731224099199
$ws.Range("S4").Value = 0.190002651698962
$ws.Range("T4").Value = 0.1900026516989619

# Row 5
$ws.Range("M5").Value = 6.457974333333333
$ws.Range("N5").Value = 19.373923
$ws.Range("O5").Value = 0.1373888468646722
$ws.Range("P5").Value = 0.1373888468646721
$ws.Range("Q5").Value = 5.673623213336445
$ws.Range("R5").Value = 51.062608920028
$ws.Range("S5").Value = 0.1373888468646722
$ws.Range("T5").Value = 0.1373888468646721

# Row 6
$ws.Range("M6").Value = 1.948535
$ws.Range("N6").Value = 5.845605
$ws.Range("O6").Value = 0.04145370713904261
$ws.Range("P6").Value = 0.0414537071390426
$ws.Range("Q6").Value = 1.711876331086666
$ws.Range("R6").Value = 15.40688697978
$ws.Range("S6").Value = 0.04145370713904261
$ws.Range("T6").Value = 0.0414537071390426

# Row 7
$ws.Range("M7").Value = 4.881814666666666
$ws.Range("N7").Value = 14.645444
$ws.Range("O7").Value = 0.1038571621752152
$ws.Range("P7").Value = 0.1038571621752152
$ws.Range("Q7").Value = 4.288895493598222
$ws.Range("R7").Value = 38.60005944238399
$ws.Range("S7").Value = 0.1038571621752152
$ws.Range("T7").Value = 0.1038571621752152
